$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 262.66
$ws.Range("J17").Value = 262.66
$ws.Range("L17").Value = 787.98
$ws.Range("N17").Value = -1123.98
$ws.Range("H51").Value = 16063.625
$ws.Range("I51").Value = 52000
$ws.Range("J51").Value = 4084.8333
$ws.Range("K51").Value = 52000
$ws.Range("L51").Value = 4084.8333
$ws.Range("M51").Value = -51516
$ws.Range("N51").Value = -5052.8333
$ws.Range("H86").Value = 4912.9
$ws.Range("I86").Value = 1419.625
$ws.Range("K86").Value = 1419.625
$ws.Range("M86").Value = -296.625
$ws.Range("H89").Value = 4912.9
$ws.Range("I89").Value = 1419.625
$ws.Range("K89").Value = 7098.125
$ws.Range("M89").Value = -1482.125
$ws.Range("H92").Value = 415.6
$ws.Range("I92").Value = 415.6
$ws.Range("K92").Value = 415.6
$ws.Range("M92").Value = 832.4
$ws.Range("H125").Value = 2133.75
$ws.Range("I125").Value = 3006.4
$ws.Range("J125").Value = 1737.091
$ws.Range("K125").Value = 27057.6
$ws.Range("L125").Value = 15633.819
$ws.Range("M125").Value = -24597.6
$ws.Range("N125").Value = -20553.819
$ws.Range("H132").Value = 5439621.5
$ws.Range("I132").Value = 6255442
$ws.Range("J132").Value = 817.3333
$ws.Range("K132").Value = 18766326
$ws.Range("L132").Value = 2451.9999
$ws.Range("M132").Value = -18763796
$ws.Range("N132").Value = -7511.9999
$ws.Range("H137").Value = 1012.5714
$ws.Range("I137").Value = 1036.1765
$ws.Range("J137").Value = 912.25
$ws.Range("K137").Value = 3108.5295
$ws.Range("L137").Value = 2736.75
$ws.Range("M137").Value = -558.5295000000001
$ws.Range("N137").Value = -7836.75
$ws.Range("H138").Value = 4219.8247
$ws.Range("I138").Value = 2172.0667
$ws.Range("J138").Value = 4951.1665
$ws.Range("K138").Value = 6516.2001
$ws.Range("L138").Value = 14853.4995
$ws.Range("M138").Value = -1376.2001
$ws.Range("N138").Value = -25133.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 64310.18
$ws.Range("I32").Value = 12079.571
$ws.Range("K32").Value = 12079.571
$ws.Range("M32").Value = -11792.571
$ws.Range("H45").Value = 1396
$ws.Range("I45").Value = 1302.4546
$ws.Range("J45").Value = 1601.8
$ws.Range("K45").Value = 1302.4546
$ws.Range("L45").Value = 1601.8
$ws.Range("M45").Value = -925.4546
$ws.Range("N45").Value = -2355.8
$ws.Range("H110").Value = 50055630
$ws.Range("I110").Value = 55617188
$ws.Range("J110").Value = 1644.5
$ws.Range("K110").Value = 55617188
$ws.Range("L110").Value = 1644.5
$ws.Range("M110").Value = -55615143
$ws.Range("N110").Value = -5734.5
$ws.Range("H112").Value = 18875.25
$ws.Range("J112").Value = 18875.25
$ws.Range("L112").Value = 18875.25
$ws.Range("N112").Value = -21829.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 106473.1
$ws.Range("I105").Value = 92130
$ws.Range("J105").Value = 126194.875
$ws.Range("K105").Value = 92130
$ws.Range("L105").Value = 126194.875
$ws.Range("M105").Value = -90383
$ws.Range("N105").Value = -129688.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 11499.333
$ws.Range("J80").Value = 11499.333
$ws.Range("L80").Value = 11499.333
$ws.Range("N80").Value = -13745.333
$ws.Range("H83").Value = 11499.333
$ws.Range("J83").Value = 11499.333
$ws.Range("L83").Value = 34497.999
$ws.Range("N83").Value = -45729.999
$ws.Range("H86").Value = 1941.091
$ws.Range("I86").Value = 1667.5
$ws.Range("J86").Value = 2269.4
$ws.Range("K86").Value = 1667.5
$ws.Range("L86").Value = 2269.4
$ws.Range("M86").Value = -544.5
$ws.Range("N86").Value = -4515.4
$ws.Range("H89").Value = 1941.091
$ws.Range("I89").Value = 1667.5
$ws.Range("J89").Value = 2269.4
$ws.Range("K89").Value = 8337.5
$ws.Range("L89").Value = 11347
$ws.Range("M89").Value = -2721.5
$ws.Range("N89").Value = -22579

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5577.156
$ws.Range("I5").Value = 999.37933
$ws.Range("J5").Value = 13874.375
$ws.Range("K5").Value = 2998.13799
$ws.Range("L5").Value = 41623.125
$ws.Range("M5").Value = -2886.13799
$ws.Range("N5").Value = -41847.125
$ws.Range("H23").Value = 572.58826
$ws.Range("I23").Value = 386.66666
$ws.Range("J23").Value = 612.4286
$ws.Range("K23").Value = 1159.99998
$ws.Range("L23").Value = 1837.2858
$ws.Range("M23").Value = -924.9999800000001
$ws.Range("N23").Value = -2307.2858
$ws.Range("H47").Value = 102.6
$ws.Range("I47").Value = 102.6
$ws.Range("K47").Value = 307.8
$ws.Range("M47").Value = 123.2
$ws.Range("H131").Value = 813.03
$ws.Range("J131").Value = 826.086
$ws.Range("L131").Value = 2478.258
$ws.Range("N131").Value = -12558.258
$ws.Range("H135").Value = 5577.156
$ws.Range("I135").Value = 999.37933
$ws.Range("J135").Value = 13874.375
$ws.Range("K135").Value = 8994.41397
$ws.Range("L135").Value = 124869.375
$ws.Range("M135").Value = -6459.41397
$ws.Range("N135").Value = -129939.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 57372.79
$ws.Range("I70").Value = 91403.44
$ws.Range("J70").Value = 5192.467
$ws.Range("K70").Value = 91403.44
$ws.Range("L70").Value = 5192.467
$ws.Range("M70").Value = -91133.44
$ws.Range("N70").Value = -5732.467
$ws.Range("H73").Value = 57372.79
$ws.Range("I73").Value = 91403.44
$ws.Range("J73").Value = 5192.467
$ws.Range("K73").Value = 91403.44
$ws.Range("L73").Value = 5192.467
$ws.Range("M73").Value = -90467.44
$ws.Range("N73").Value = -7064.467
$ws.Range("H122").Value = 5208.1816
$ws.Range("J122").Value = 4470
$ws.Range("L122").Value = 13410
$ws.Range("N122").Value = -18310

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1999.32
$ws.Range("I7").Value = 1615.9474
$ws.Range("K7").Value = 1615.9474
$ws.Range("M7").Value = -1503.9474
$ws.Range("H61").Value = 2798.7334
$ws.Range("I61").Value = 2203.1667
$ws.Range("K61").Value = 2203.1667
$ws.Range("M61").Value = -2001.1667
$ws.Range("H113").Value = 2798.7334
$ws.Range("I113").Value = 2203.1667
$ws.Range("K113").Value = 2203.1667
$ws.Range("M113").Value = -33.16670000000022
$ws.Range("H126").Value = 1999.32
$ws.Range("I126").Value = 1615.9474
$ws.Range("K126").Value = 4847.8422
$ws.Range("M126").Value = -2377.8422

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 444.79166
$ws.Range("I113").Value = 338.66666
$ws.Range("K113").Value = 1015.99998
$ws.Range("M113").Value = 1154.00002
$ws.Range("H132").Value = 3447.9546
$ws.Range("I132").Value = 3966
$ws.Range("J132").Value = 2699.6667
$ws.Range("K132").Value = 11898
$ws.Range("L132").Value = 8099.000100000001
$ws.Range("M132").Value = -9368
$ws.Range("N132").Value = -13159.0001

Write-Host "Applied all cell updates"
